$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '27.559.44'
Set-TextValue $ws.Cells.Item(2, 5) '  +0.56%  '

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '1.639.28'
Set-TextValue $ws.Cells.Item(3, 5) '  -0.75%  '

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) '0.999'
Set-TextValue $ws.Cells.Item(4, 5) '  -0.03%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '212.64'
Set-TextValue $ws.Cells.Item(5, 5) '  -0.40%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '0.533'
Set-TextValue $ws.Cells.Item(6, 5) '  +4.65%  '

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) '0.998'
Set-TextValue $ws.Cells.Item(7, 5) '  -0.11%  '

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) '23.05'
Set-TextValue $ws.Cells.Item(8, 5) '  -4.89%  '

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '0.257'
Set-TextValue $ws.Cells.Item(9, 5) '  -2.09%  '

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) '0.0611'
Set-TextValue $ws.Cells.Item(10, 5) '  -0.68%  '

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.0889'
Set-TextValue $ws.Cells.Item(11, 5) '  +1.32%  '

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) '1.867.49'
Set-TextValue $ws.Cells.Item(12, 5) '  -0.90%  '

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) '1.631.34'
Set-TextValue $ws.Cells.Item(13, 5) '  -1.54%  '

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) '0.566'
Set-TextValue $ws.Cells.Item(14, 5) '  -1.02%  '

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) '4.03'
Set-TextValue $ws.Cells.Item(15, 5) '  -1.28%  '

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '64.38'
Set-TextValue $ws.Cells.Item(16, 5) '  -2.39%  '

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) '27.465.49'
Set-TextValue $ws.Cells.Item(17, 5) '  +0.26%  '

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) '230.35'
Set-TextValue $ws.Cells.Item(18, 5) '  -1.78%  '

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) '7.70'
Set-TextValue $ws.Cells.Item(19, 5) '  +3.06%  '

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) '0.0₃0725'
Set-TextValue $ws.Cells.Item(20, 5) '  -0.04%  '

# Row 21
Set-TextValue $ws.Cells.Item(21, 5) '  +0.13%  '

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) '4.32'
Set-TextValue $ws.Cells.Item(22, 5) '  -1.78%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '10.00'
Set-TextValue $ws.Cells.Item(23, 5) '  +7.48%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '1.95'
Set-TextValue $ws.Cells.Item(24, 5) '  -3.39%  '

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '149.50'
Set-TextValue $ws.Cells.Item(25, 5) '  +1.72%  '

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) '6.98'
Set-TextValue $ws.Cells.Item(26, 5) '  -3.05%  '

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) '0.113'
Set-TextValue $ws.Cells.Item(27, 5) '  +1.77%  '

# Row 28
Set-TextValue $ws.Cells.Item(28, 5) '  -0.06%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '15.61'
Set-TextValue $ws.Cells.Item(29, 5) '  -2.98%  '

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) '1.19'
Set-TextValue $ws.Cells.Item(30, 5) '  -0.56%  '

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) '0.0488'
Set-TextValue $ws.Cells.Item(31, 5) '  -1.95%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) '3.30'
Set-TextValue $ws.Cells.Item(32, 5) '  -0.44%  '

# Row 33
Set-TextValue $ws.Cells.Item(33, 5) '  +2.21%  '

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) '1.415.38'
Set-TextValue $ws.Cells.Item(34, 5) '  -2.89%  '

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) '1.59'
Set-TextValue $ws.Cells.Item(35, 5) '  +2.62%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 5) '  -1.88%  '

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) '0.572'
Set-TextValue $ws.Cells.Item(37, 5) '  -0.26%  '

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) '0.878'
Set-TextValue $ws.Cells.Item(38, 5) '  -3.50%  '

# Row 39
Set-TextValue $ws.Cells.Item(39, 5) '  -1.61%  '

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '0.896'
Set-TextValue $ws.Cells.Item(40, 5) '  +14.01%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 5) '  -0.37%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 5) '  +0.06%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '2.47'
Set-TextValue $ws.Cells.Item(43, 5) '  -0.50%  '

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) '5.51'
Set-TextValue $ws.Cells.Item(44, 5) '  +1.38%  '

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) '2.25'
Set-TextValue $ws.Cells.Item(45, 5) '  +1.58%  '

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '65.03'
Set-TextValue $ws.Cells.Item(46, 5) '  -0.49%  '

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '1.777.69'
Set-TextValue $ws.Cells.Item(47, 5) '  -0.86%  '

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) '1.68'
Set-TextValue $ws.Cells.Item(48, 5) '  -2.72%  '

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) '86.19'
Set-TextValue $ws.Cells.Item(49, 5) '  -2.63%  '

# Row 50
Set-TextValue $ws.Cells.Item(50, 5) '  +0.38%  '

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) '0.0991'
Set-TextValue $ws.Cells.Item(51, 5) '  -1.95%  '
